$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 236, shifting existing rows 236:292 down to 237:293.
$ws.Rows(236).Insert()

# Populate the newly inserted row 236 with the new record.
$ws.Cells.Item(236, 1).Value = 8
$ws.Cells.Item(236, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(236, 3).Value = "Coquimbo"
$ws.Cells.Item(236, 4).Value = 44754
$ws.Cells.Item(236, 5).Value = 4
$ws.Cells.Item(236, 6).Value = 100112003
$ws.Cells.Item(236, 7).Value = "Ajo"
$ws.Cells.Item(236, 8).Value = "Chino"
$ws.Cells.Item(236, 9).Value = "Primera"
$ws.Cells.Item(236, 10).Value = 520
$ws.Cells.Item(236, 11).Value = 24000
$ws.Cells.Item(236, 12).Value = 25000
$ws.Cells.Item(236, 13).Value = 24500
$ws.Cells.Item(236, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(236, 15).Value = "China"
$ws.Cells.Item(236, 16).Value = 2450
$ws.Cells.Item(236, 17).Value = 10
$ws.Cells.Item(236, 18).Value = "Hortaliza"
